# Sync attendance_reports: reorder "System, <email>" -> "<email>, System"
# in column G for the rows where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70, `
          86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142, `
          144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248, `
          251,255,272,275,278,282,299,302,305,309)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G is the 7th column
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
